$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("missing_stations")

# --- Restore original (pre-sort) row order for the data rows (A2:G58) ---
# The sheet had previously been sorted (ascending) by column C (free/taken),
# this undoes that sort, putting rows back in their original entry order,
# while preserving each row's associated E/F/G annotations.

$ws.Cells.Item(2,1).Value = 651
$ws.Cells.Item(2,2).Value = "Michigan Ave & 71st St"
$ws.Cells.Item(2,5).ClearContents()
$ws.Cells.Item(2,6).ClearContents()
$ws.Cells.Item(2,7).ClearContents()

$ws.Cells.Item(3,1).Value = 459
$ws.Cells.Item(3,2).Value = "Lakefront Trail & Bryn Mawr Ave"
$ws.Cells.Item(3,5).ClearContents()
$ws.Cells.Item(3,6).ClearContents()
$ws.Cells.Item(3,7).ClearContents()

$ws.Cells.Item(4,1).Value = 20
$ws.Cells.Item(4,2).Value = "Damen Ave & Wabansia Ave"
$ws.Cells.Item(4,5).ClearContents()
$ws.Cells.Item(4,6).ClearContents()
$ws.Cells.Item(4,7).ClearContents()

$ws.Cells.Item(5,1).Value = 357
$ws.Cells.Item(5,2).Value = "Lamon Ave & Armitage Ave"
$ws.Cells.Item(5,5).ClearContents()
$ws.Cells.Item(5,6).ClearContents()
$ws.Cells.Item(5,7).ClearContents()

$ws.Cells.Item(6,1).Value = 358
$ws.Cells.Item(6,2).Value = "Kilpatrick Ave & Parker Ave"
$ws.Cells.Item(6,5).ClearContents()
$ws.Cells.Item(6,6).ClearContents()
$ws.Cells.Item(6,7).ClearContents()

$ws.Cells.Item(7,1).Value = 360
$ws.Cells.Item(7,2).Value = "Kilbourn & Roscoe"
$ws.Cells.Item(7,5).ClearContents()
$ws.Cells.Item(7,6).ClearContents()
$ws.Cells.Item(7,7).Value = "confirmed"

$ws.Cells.Item(8,1).Value = 361
$ws.Cells.Item(8,2).Value = "Kenosha & Wellington"
$ws.Cells.Item(8,5).ClearContents()
$ws.Cells.Item(8,6).ClearContents()
$ws.Cells.Item(8,7).Value = "confirmed"

$ws.Cells.Item(9,1).Value = 362
$ws.Cells.Item(9,2).Value = "Lawndale Ave & 16th St"
$ws.Cells.Item(9,5).ClearContents()
$ws.Cells.Item(9,6).ClearContents()
$ws.Cells.Item(9,7).ClearContents()

$ws.Cells.Item(10,1).Value = 363
$ws.Cells.Item(10,2).Value = "Tripp Ave & 15th St"
$ws.Cells.Item(10,5).ClearContents()
$ws.Cells.Item(10,6).ClearContents()
$ws.Cells.Item(10,7).Value = "confirmed"

$ws.Cells.Item(11,1).Value = 473
$ws.Cells.Item(11,2).Value = "Mason Ave & Belmont Ave"
$ws.Cells.Item(11,5).ClearContents()
$ws.Cells.Item(11,6).ClearContents()
$ws.Cells.Item(11,7).ClearContents()

$ws.Cells.Item(12,1).Value = 329
$ws.Cells.Item(12,2).Value = "Central Park Ave & Douglas Blvd"
$ws.Cells.Item(12,5).ClearContents()
$ws.Cells.Item(12,6).ClearContents()
$ws.Cells.Item(12,7).ClearContents()

$ws.Cells.Item(13,1).Value = 330
$ws.Cells.Item(13,2).Value = "Keeler Ave & Roosevelt Rd"
$ws.Cells.Item(13,5).ClearContents()
$ws.Cells.Item(13,6).ClearContents()
$ws.Cells.Item(13,7).ClearContents()

$ws.Cells.Item(14,1).Value = 331
$ws.Cells.Item(14,2).Value = "Pulaski Rd & 21st St"
$ws.Cells.Item(14,5).ClearContents()
$ws.Cells.Item(14,6).ClearContents()
$ws.Cells.Item(14,7).ClearContents()

$ws.Cells.Item(15,1).Value = 332
$ws.Cells.Item(15,2).Value = "Harding Ave & 26th St"
$ws.Cells.Item(15,5).ClearContents()
$ws.Cells.Item(15,6).ClearContents()
$ws.Cells.Item(15,7).ClearContents()

$ws.Cells.Item(16,1).Value = 334
$ws.Cells.Item(16,2).Value = "Lawndale Ave & 30th St"
$ws.Cells.Item(16,5).ClearContents()
$ws.Cells.Item(16,6).ClearContents()
$ws.Cells.Item(16,7).ClearContents()

$ws.Cells.Item(17,1).Value = 335
$ws.Cells.Item(17,2).Value = "Komensky Ave & 31st St"
$ws.Cells.Item(17,5).ClearContents()
$ws.Cells.Item(17,6).ClearContents()
$ws.Cells.Item(17,7).ClearContents()

$ws.Cells.Item(18,1).Value = 364
$ws.Cells.Item(18,2).Value = "Homan Ave & Fillmore St"
$ws.Cells.Item(18,5).ClearContents()
$ws.Cells.Item(18,6).ClearContents()
$ws.Cells.Item(18,7).ClearContents()

$ws.Cells.Item(19,1).Value = 365
$ws.Cells.Item(19,2).Value = "Kildare Ave & 26th St"
$ws.Cells.Item(19,5).ClearContents()
$ws.Cells.Item(19,6).ClearContents()
$ws.Cells.Item(19,7).ClearContents()

$ws.Cells.Item(20,1).Value = 366
$ws.Cells.Item(20,2).Value = "Spaulding Ave & 16th St"
$ws.Cells.Item(20,5).ClearContents()
$ws.Cells.Item(20,6).ClearContents()
$ws.Cells.Item(20,7).ClearContents()

$ws.Cells.Item(21,1).Value = 368
$ws.Cells.Item(21,2).Value = "Tripp Ave & 31st St"
$ws.Cells.Item(21,5).ClearContents()
$ws.Cells.Item(21,6).ClearContents()
$ws.Cells.Item(21,7).ClearContents()

$ws.Cells.Item(22,1).Value = 397
$ws.Cells.Item(22,2).Value = "Narragansett Ave & School St"
$ws.Cells.Item(22,5).ClearContents()
$ws.Cells.Item(22,6).ClearContents()
$ws.Cells.Item(22,7).Value = "confirmed"

$ws.Cells.Item(23,1).Value = 371
$ws.Cells.Item(23,2).Value = "Kildare Ave & Chicago Ave"
$ws.Cells.Item(23,5).ClearContents()
$ws.Cells.Item(23,6).ClearContents()
$ws.Cells.Item(23,7).ClearContents()

$ws.Cells.Item(24,1).Value = 379
$ws.Cells.Item(24,2).Value = "Rockwell St & Archer Ave"
$ws.Cells.Item(24,5).ClearContents()
$ws.Cells.Item(24,6).ClearContents()
$ws.Cells.Item(24,7).ClearContents()

$ws.Cells.Item(25,1).Value = 380
$ws.Cells.Item(25,2).Value = "Fairfield Ave & 44th St"
$ws.Cells.Item(25,5).ClearContents()
$ws.Cells.Item(25,6).ClearContents()
$ws.Cells.Item(25,7).ClearContents()

$ws.Cells.Item(26,1).Value = 387
$ws.Cells.Item(26,2).Value = "St Louis Ave & 59th St"
$ws.Cells.Item(26,5).ClearContents()
$ws.Cells.Item(26,6).ClearContents()
$ws.Cells.Item(26,7).ClearContents()

$ws.Cells.Item(27,1).Value = 389
$ws.Cells.Item(27,2).Value = "Maplewood Ave & 59th St"
$ws.Cells.Item(27,5).ClearContents()
$ws.Cells.Item(27,6).ClearContents()
$ws.Cells.Item(27,7).ClearContents()

$ws.Cells.Item(28,1).Value = 404
$ws.Cells.Item(28,2).Value = "Pulaski Rd & 60th St"
$ws.Cells.Item(28,5).ClearContents()
$ws.Cells.Item(28,6).ClearContents()
$ws.Cells.Item(28,7).ClearContents()

$ws.Cells.Item(29,1).Value = 1
$ws.Cells.Item(29,2).Value = "Special Events"
$ws.Cells.Item(29,5).ClearContents()
$ws.Cells.Item(29,6).ClearContents()
$ws.Cells.Item(29,7).ClearContents()

$ws.Cells.Item(30,1).Value = 372
$ws.Cells.Item(30,2).Value = "California Ave & Augusta Blvd"
$ws.Cells.Item(30,5).Value = "California Ave & Cortez St"
$ws.Cells.Item(30,6).Value = 622
$ws.Cells.Item(30,7).ClearContents()

$ws.Cells.Item(31,1).Value = 372
$ws.Cells.Item(31,2).Value = "Humboldt Dr & Luis Munoz Marin Dr"
$ws.Cells.Item(31,5).ClearContents()
$ws.Cells.Item(31,6).ClearContents()
$ws.Cells.Item(31,7).ClearContents()

$ws.Cells.Item(32,1).Value = 397
$ws.Cells.Item(32,2).Value = "Saginaw Ave & Exchange Ave"
$ws.Cells.Item(32,5).Value = "Kingston Ave & 75th St"
$ws.Cells.Item(32,6).ClearContents()
$ws.Cells.Item(32,7).ClearContents()

$ws.Cells.Item(33,1).Value = 512
$ws.Cells.Item(33,2).Value = "BBB ~ Divvy Parts Testing"
$ws.Cells.Item(33,5).ClearContents()
$ws.Cells.Item(33,6).ClearContents()
$ws.Cells.Item(33,7).ClearContents()

$ws.Cells.Item(34,1).Value = 606
$ws.Cells.Item(34,2).Value = "Forest Ave & Chicago Ave"
$ws.Cells.Item(34,5).ClearContents()
$ws.Cells.Item(34,6).ClearContents()
$ws.Cells.Item(34,7).ClearContents()

$ws.Cells.Item(35,1).Value = 607
$ws.Cells.Item(35,2).Value = "Cuyler Ave & Augusta St"
$ws.Cells.Item(35,5).ClearContents()
$ws.Cells.Item(35,6).ClearContents()
$ws.Cells.Item(35,7).ClearContents()

$ws.Cells.Item(36,1).Value = 608
$ws.Cells.Item(36,2).Value = "Humphrey Ave & Ontario St"
$ws.Cells.Item(36,5).ClearContents()
$ws.Cells.Item(36,6).ClearContents()
$ws.Cells.Item(36,7).ClearContents()

$ws.Cells.Item(37,1).Value = 609
$ws.Cells.Item(37,2).Value = "Forest Ave & Lake St"
$ws.Cells.Item(37,5).ClearContents()
$ws.Cells.Item(37,6).ClearContents()
$ws.Cells.Item(37,7).ClearContents()

$ws.Cells.Item(38,1).Value = 610
$ws.Cells.Item(38,2).Value = "Marion St & South Blvd"
$ws.Cells.Item(38,5).ClearContents()
$ws.Cells.Item(38,6).ClearContents()
$ws.Cells.Item(38,7).ClearContents()

$ws.Cells.Item(39,1).Value = 611
$ws.Cells.Item(39,2).Value = "Oak Park Ave & South Blvd"
$ws.Cells.Item(39,5).ClearContents()
$ws.Cells.Item(39,6).ClearContents()
$ws.Cells.Item(39,7).ClearContents()

$ws.Cells.Item(40,1).Value = 612
$ws.Cells.Item(40,2).Value = "Ridgeland Ave & Lake St"
$ws.Cells.Item(40,5).ClearContents()
$ws.Cells.Item(40,6).ClearContents()
$ws.Cells.Item(40,7).ClearContents()

$ws.Cells.Item(41,1).Value = 613
$ws.Cells.Item(41,2).Value = "Wisconsin Ave & Madison St"
$ws.Cells.Item(41,5).ClearContents()
$ws.Cells.Item(41,6).ClearContents()
$ws.Cells.Item(41,7).ClearContents()

$ws.Cells.Item(42,1).Value = 613
$ws.Cells.Item(42,2).Value = "Wisconsin Ave & Madison St (Temp)"
$ws.Cells.Item(42,5).Value = "Wisconsin Ave & Madison St"
$ws.Cells.Item(42,6).ClearContents()
$ws.Cells.Item(42,7).ClearContents()

$ws.Cells.Item(43,1).Value = 614
$ws.Cells.Item(43,2).Value = "East Ave & Madison St"
$ws.Cells.Item(43,5).ClearContents()
$ws.Cells.Item(43,6).ClearContents()
$ws.Cells.Item(43,7).ClearContents()

$ws.Cells.Item(44,1).Value = 615
$ws.Cells.Item(44,2).Value = "Lombard Ave & Madison St"
$ws.Cells.Item(44,5).ClearContents()
$ws.Cells.Item(44,6).ClearContents()
$ws.Cells.Item(44,7).ClearContents()

$ws.Cells.Item(45,1).Value = 616
$ws.Cells.Item(45,2).Value = "Oak Park Ave & Harrison St"
$ws.Cells.Item(45,5).ClearContents()
$ws.Cells.Item(45,6).ClearContents()
$ws.Cells.Item(45,7).ClearContents()

$ws.Cells.Item(46,1).Value = 617
$ws.Cells.Item(46,2).Value = "East Ave & Garfield St"
$ws.Cells.Item(46,5).ClearContents()
$ws.Cells.Item(46,6).ClearContents()
$ws.Cells.Item(46,7).ClearContents()

$ws.Cells.Item(47,1).Value = 618
$ws.Cells.Item(47,2).Value = "Lombard Ave & Garfield St"
$ws.Cells.Item(47,5).ClearContents()
$ws.Cells.Item(47,6).ClearContents()
$ws.Cells.Item(47,7).ClearContents()

$ws.Cells.Item(48,1).Value = 669
$ws.Cells.Item(48,2).Value = "LBS - BBB La Magie"
$ws.Cells.Item(48,5).ClearContents()
$ws.Cells.Item(48,6).ClearContents()
$ws.Cells.Item(48,7).ClearContents()

$ws.Cells.Item(49,1).Value = 670
$ws.Cells.Item(49,2).Value = "MTL-ECO5.1-01"
$ws.Cells.Item(49,5).ClearContents()
$ws.Cells.Item(49,6).ClearContents()
$ws.Cells.Item(49,7).ClearContents()

$ws.Cells.Item(50,1).Value = 671
$ws.Cells.Item(50,2).Value = "HUBBARD ST BIKE CHECKING (LBS-WH-TEST)"
$ws.Cells.Item(50,5).Value = "Base - 2132 W Hubbard Warehouse"
$ws.Cells.Item(50,6).ClearContents()
$ws.Cells.Item(50,7).ClearContents()

$ws.Cells.Item(51,1).Value = 675
$ws.Cells.Item(51,2).Value = "HQ QR"
$ws.Cells.Item(51,5).ClearContents()
$ws.Cells.Item(51,6).ClearContents()
$ws.Cells.Item(51,7).ClearContents()

$ws.Cells.Item(52,1).Value = 676
$ws.Cells.Item(52,2).Value = "WATSON TESTING - DIVVY"
$ws.Cells.Item(52,5).ClearContents()
$ws.Cells.Item(52,6).ClearContents()
$ws.Cells.Item(52,7).ClearContents()

$ws.Cells.Item(53,1).ClearContents()
$ws.Cells.Item(53,2).Value = "DIVVY CASSETTE REPAIR MOBILE STATION"
$ws.Cells.Item(53,5).ClearContents()
$ws.Cells.Item(53,6).ClearContents()
$ws.Cells.Item(53,7).ClearContents()

$ws.Cells.Item(54,1).ClearContents()
$ws.Cells.Item(54,2).Value = "DIVVY Map Frame B/C Station"
$ws.Cells.Item(54,5).ClearContents()
$ws.Cells.Item(54,6).ClearContents()
$ws.Cells.Item(54,7).ClearContents()

$ws.Cells.Item(55,1).ClearContents()
$ws.Cells.Item(55,2).Value = "TS ~ DIVVY PARTS TESTING"
$ws.Cells.Item(55,5).ClearContents()
$ws.Cells.Item(55,6).ClearContents()
$ws.Cells.Item(55,7).ClearContents()

$ws.Cells.Item(56,1).Value = 566
$ws.Cells.Item(56,2).Value = "Ashland Ave & 69th St"
$ws.Cells.Item(56,5).ClearContents()
$ws.Cells.Item(56,6).ClearContents()
$ws.Cells.Item(56,7).ClearContents()

$ws.Cells.Item(57,1).Value = 625
$ws.Cells.Item(57,2).Value = "Chicago Ave & Dempster St"
$ws.Cells.Item(57,5).ClearContents()
$ws.Cells.Item(57,6).ClearContents()
$ws.Cells.Item(57,7).ClearContents()

$ws.Cells.Item(58,1).Value = 704
$ws.Cells.Item(58,2).Value = "Jeffery Blvd & 91st St"
$ws.Cells.Item(58,5).ClearContents()
$ws.Cells.Item(58,6).ClearContents()
$ws.Cells.Item(58,7).ClearContents()

# --- Remove the stale AutoFilter sort-state (the column-C sort no longer applies) ---
$ws.AutoFilterMode = $false
$ws.Range("A1:G58").AutoFilter() | Out-Null

# --- Best-effort restore of the view / frozen-pane scroll position and selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D27").Select()

$wb.Save()